# Update the LR-pairs (Hspg2-Itgb1) sheet with newly recomputed TPM-based
# expression / specificity / edge-weight values.
#
# Only numeric values change (ligand/receptor average & total expression,
# derived specificity scores, and edge weights in columns G,H,I,J,M,N,O,P,Q,R,S,T);
# the categorical columns (A-F,K,L) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 91.34108066666666
$ws.Range("H2").Value = 274.023242
$ws.Range("I2").Value = 0.2190334467302001
$ws.Range("J2").Value = 0.2190334467302
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 5575.606288115911
$ws.Range("R2").Value = 50180.4565930432
$ws.Range("S2").Value = 0.0447619700053872
$ws.Range("T2").Value = 0.0447619700053872

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 91.34108066666666
$ws.Range("H3").Value = 274.023242
$ws.Range("I3").Value = 0.2190334467302001
$ws.Range("J3").Value = 0.2190334467302
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 9710.878214939588
$ws.Range("R3").Value = 87397.90393445631
$ws.Range("S3").Value = 0.07796067672668809
$ws.Range("T3").Value = 0.07796067672668809

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 91.34108066666666
$ws.Range("H4").Value = 274.023242
$ws.Range("I4").Value = 0.2190334467302001
$ws.Range("J4").Value = 0.2190334467302
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 11996.59224154768
$ws.Range("R4").Value = 107969.3301739292
$ws.Range("S4").Value = 0.09631079999812475
$ws.Range("T4").Value = 0.09631079999812475

# Row 5 (FAPs -> ECs)
$ws.Range("G5").Value = 276.4348856666666
$ws.Range("H5").Value = 829.3046569999999
$ws.Range("I5").Value = 0.6628833966285105
$ws.Range("J5").Value = 0.6628833966285105
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 16874.02946766467
$ws.Range("R5").Value = 151866.2652089821
$ws.Range("S5").Value = 0.1354677432141392
$ws.Range("T5").Value = 0.1354677432141392

# Row 6 (FAPs -> FAPs)
$ws.Range("G6").Value = 276.4348856666666
$ws.Range("H6").Value = 829.3046569999999
$ws.Range("I6").Value = 0.6628833966285105
$ws.Range("J6").Value = 0.6628833966285105
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("Q6").Value = 29389.02725342272
$ws.Range("S6").Value = 0.2359403961519218
$ws.Range("T6").Value = 0.2359403961519219

# Row 7 (FAPs -> MuSCs)
$ws.Range("G7").Value = 276.4348856666666
$ws.Range("H7").Value = 829.3046569999999
$ws.Range("I7").Value = 0.6628833966285105
$ws.Range("J7").Value = 0.6628833966285105
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 36306.51816770186
$ws.Range("R7").Value = 326758.6635093167
$ws.Range("S7").Value = 0.2914752572624494
$ws.Range("T7").Value = 0.2914752572624494

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 49.24290466666667
$ws.Range("H8").Value = 147.728714
$ws.Range("I8").Value = 0.1180831566412894
$ws.Range("J8").Value = 0.1180831566412894
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 3005.86600137253
$ws.Range("R8").Value = 27052.79401235277
$ws.Range("S8").Value = 0.02413163283792703
$ws.Range("T8").Value = 0.02413163283792703

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 49.24290466666667
$ws.Range("H9").Value = 147.728714
$ws.Range("I9").Value = 0.1180831566412894
$ws.Range("J9").Value = 0.1180831566412894
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 5235.233113925575
$ws.Range("R9").Value = 47117.09802533017
$ws.Range("S9").Value = 0.04202939294982636
$ws.Range("T9").Value = 0.04202939294982636

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 49.24290466666667
$ws.Range("H10").Value = 147.728714
$ws.Range("I10").Value = 0.1180831566412894
$ws.Range("J10").Value = 0.1180831566412894
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 6467.484769872976
$ws.Range("R10").Value = 58207.36292885679
$ws.Range("S10").Value = 0.05192213085353604
$ws.Range("T10").Value = 0.05192213085353604
